{"js": "// Replace each multiplication expression's old value with the new one.\n// Each (old, new) pair is unique in this document, so a direct search+\n// replace for each pair is unambiguous and order-independent.\nconst replacements = [\n  [\"320\u00d72=\", \"588\u00d78=\"],\n  [\"804\u00d75=\", \"836\u00d76=\"],\n  [\"308\u00d77=\", \"499\u00d78=\"],\n  [\"721\u00d73=\", \"151\u00d73=\"],\n  [\"526\u00d74=\", \"590\u00d75=\"],\n  [\"948\u00d79=\", \"401\u00d79=\"],\n  [\"829\u00d73=\", \"124\u00d78=\"],\n  [\"656\u00d72=\", \"236\u00d75=\"],\n  [\"576\u00d72=\", \"311\u00d72=\"],\n  [\"214\u00d77=\", \"584\u00d77=\"],\n  [\"299\u00d75=\", \"486\u00d73=\"],\n  [\"424\u00d74=\", \"479\u00d76=\"],\n  [\"314\u00d75=\", \"359\u00d75=\"],\n  [\"662\u00d76=\", \"607\u00d76=\"],\n  [\"823\u00d78=\", \"817\u00d72=\"],\n  [\"423\u00d79=\", \"201\u00d76=\"],\n  [\"538\u00d76=\", \"248\u00d76=\"],\n  [\"748\u00d72=\", \"802\u00d79=\"],\n  [\"863\u00d75=\", \"794\u00d73=\"],\n  [\"109\u00d76=\", \"143\u00d76=\"],\n  [\"592\u00d77=\", \"252\u00d72=\"],\n  [\"365\u00d77=\", \"823\u00d77=\"],\n  [\"605\u00d76=\", \"510\u00d75=\"],\n  [\"895\u00d76=\", \"733\u00d72=\"],\n  [\"162\u00d74=\", \"393\u00d76=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update each multiplication expression's old value to the new one.\n# Each (old, new) pair is unique within the document, so Find/Replace\n# for each pair independently (order-independent) is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"320\u00d72=\", \"588\u00d78=\"),\n  @(\"804\u00d75=\", \"836\u00d76=\"),\n  @(\"308\u00d77=\", \"499\u00d78=\"),\n  @(\"721\u00d73=\", \"151\u00d73=\"),\n  @(\"526\u00d74=\", \"590\u00d75=\"),\n  @(\"948\u00d79=\", \"401\u00d79=\"),\n  @(\"829\u00d73=\", \"124\u00d78=\"),\n  @(\"656\u00d72=\", \"236\u00d75=\"),\n  @(\"576\u00d72=\", \"311\u00d72=\"),\n  @(\"214\u00d77=\", \"584\u00d77=\"),\n  @(\"299\u00d75=\", \"486\u00d73=\"),\n  @(\"424\u00d74=\", \"479\u00d76=\"),\n  @(\"314\u00d75=\", \"359\u00d75=\"),\n  @(\"662\u00d76=\", \"607\u00d76=\"),\n  @(\"823\u00d78=\", \"817\u00d72=\"),\n  @(\"423\u00d79=\", \"201\u00d76=\"),\n  @(\"538\u00d76=\", \"248\u00d76=\"),\n  @(\"748\u00d72=\", \"802\u00d79=\"),\n  @(\"863\u00d75=\", \"794\u00d73=\"),\n  @(\"109\u00d76=\", \"143\u00d76=\"),\n  @(\"592\u00d77=\", \"252\u00d72=\"),\n  @(\"365\u00d77=\", \"823\u00d77=\"),\n  @(\"605\u00d76=\", \"510\u00d75=\"),\n  @(\"895\u00d76=\", \"733\u00d72=\"),\n  @(\"162\u00d74=\", \"393\u00d76=\"),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute(\n    [ref]$oldText,   # FindText\n    [ref]$true,      # MatchCase\n    [ref]$false,     # MatchWholeWord\n    [ref]$false,     # MatchWildcards\n    [ref]$false,     # MatchSoundsLike\n    [ref]$false,     # MatchAllWordForms\n    [ref]$true,      # Forward\n    [ref]1,          # Wrap: wdFindContinue\n    [ref]$false,     # Format\n    [ref]$newText,   # ReplaceWith\n    [ref]2           # Replace: wdReplaceAll\n  )\n}\n"}
